$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.600.36'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").Value = '2.679.99'
$ws.Range("E3").Value = '  +2.04%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.37'
$ws.Range("E5").Value = '  -1.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.52'
$ws.Range("E6").Value = '  -2.23%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -0.85%  '

$ws.Range("D9").Value = '2.678.52'
$ws.Range("E9").Value = '  +2.10%  '

$ws.Range("E10").Value = '  +0.86%  '

$ws.Range("E11").Value = '  +2.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.352'
$ws.Range("E12").Value = '  +1.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.00'
$ws.Range("E13").Value = '  -1.27%  '

$ws.Range("D14").Value = '3.166.13'
$ws.Range("E14").Value = '  +0.84%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000184'
$ws.Range("E15").Value = '  -1.31%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '71.585.33'
$ws.Range("E16").Value = '  +0.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.12'
$ws.Range("E17").Value = '  -2.27%  '

$ws.Range("D18").Value = '2.670.32'
$ws.Range("E18").Value = '  +2.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.06'
$ws.Range("E19").Value = '  +5.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.00'
$ws.Range("E20").Value = '  +1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '368.14'
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").Value = '  +0.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.01'
$ws.Range("E23").Value = '  +0.72%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.58'
$ws.Range("E24").Value = '  -1.40%  '

$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.30'
$ws.Range("E26").Value = '  -2.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.84'
$ws.Range("E27").Value = '  -1.72%  '

$ws.Range("D28").Value = '2.812.80'
$ws.Range("E28").Value = '  +1.87%  '

$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("D30").Value = '0.0₃0955'
$ws.Range("E30").Value = '  -0.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.05'
$ws.Range("E31").Value = '  -0.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '505.28'
$ws.Range("E32").Value = '  -7.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("E33").Value = '  -3.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("E34").Value = '  -1.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.64'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.39'
$ws.Range("E37").Value = '  +0.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.05'
$ws.Range("E38").Value = '  -0.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.37'
$ws.Range("E39").Value = '  -2.74%  '

$ws.Range("E40").Value = '  -5.79%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.79'
$ws.Range("E41").Value = '  -3.95%  '

$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.99'
$ws.Range("E43").Value = '  -1.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.57'
$ws.Range("E44").Value = '  -1.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.332'
$ws.Range("E45").Value = '  -0.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '155.18'
$ws.Range("E46").Value = '  +1.42%  '

$ws.Range("E47").Value = '  -1.97%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.70'
$ws.Range("E48").Value = '  +1.59%  '

$ws.Range("B49").Value = 'Optimism'
$ws.Range("C49").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.73'
$ws.Range("E49").Value = '  +2.58%  '

$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.547'
$ws.Range("E50").Value = '  +1.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0761'
$ws.Range("E51").Value = '  +0.15%  '
